# Updates cryptos list prices / 1h volume percentages (and swaps the
# Toncoin / PancakeSwap rows 29-30) to match the latest scrape.
#
# Note: D-column price cells whose new text parses as a plain number
# (e.g. "258.76") would otherwise be auto-coerced to a numeric cell by
# Excel. To keep them as text (matching the original inlineStr cells)
# we briefly force a "@" (text) number format before assigning the
# value, then call ClearFormats() to drop the now-unneeded direct
# formatting so the cell's style stays the same as before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.704.38'

$ws.Range("E2").Value = '  +3.10%  '

$ws.Range("D3").Value = '2.188.82'

$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.76'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = '  +1.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.33'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").Value = '  +10.05%  '

$ws.Range("E7").Value = '  +2.76%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  +1.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.04'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").Value = '  +6.22%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0919'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("E12").Value = '  +2.71%  '

$ws.Range("E13").Value = '  +2.53%  '

$ws.Range("D14").Value = '2.517.06'

$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.28'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").Value = '  +0.72%  '

$ws.Range("D16").Value = '2.187.99'

$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.777'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").Value = '  +1.45%  '

$ws.Range("D18").Value = '43.621.12'

$ws.Range("E18").Value = '  +3.08%  '

$ws.Range("E19").Value = '  +0.85%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.07'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").Value = '  -0.75%  '

$ws.Range("E21").Value = '  +0.98%  '

$ws.Range("E22").Value = '  +12.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.32'
$ws.Range("D23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("D24").ClearFormats()

$ws.Range("E24").Value = '  -5.09%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '41.83'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").Value = '  +14.57%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.66'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = '  +1.90%  '

$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("B29").Value = 'Toncoin'

$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").Value = '  +1.23%  '

$ws.Range("B30").Value = 'PancakeSwap'

$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("D30").ClearFormats()

$ws.Range("E30").Value = '  +2.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.11'
$ws.Range("D31").ClearFormats()

$ws.Range("E31").Value = '  +1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.35'
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = '  +1.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0873'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").Value = '  +8.21%  '

$ws.Range("E34").Value = '  +3.12%  '

$ws.Range("E35").Value = '  +5.32%  '

$ws.Range("E36").Value = '  +1.40%  '

$ws.Range("E37").Value = '  +5.75%  '

$ws.Range("E38").Value = '  +4.41%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.21'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").Value = '  +11.88%  '

$ws.Range("E40").Value = '  +15.02%  '

$ws.Range("E41").Value = '  +1.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '62.83'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").Value = '  +5.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").Value = '  +6.00%  '

$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.04'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0984'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").Value = '  +1.43%  '

$ws.Range("E47").Value = '  -0.49%  '

$ws.Range("E48").Value = '  +4.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.55'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").Value = '  +27.16%  '

$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("E51").Value = '  -6.23%  '
